$wb = $excel.ActiveWorkbook

$wsStd = $wb.Worksheets.Item("IAM&CSV Standard price list")
$wsUaw = $wb.Worksheets.Item("IAM UAW East")

# ---------------------------------------------------------------------------
# "IAM&CSV Standard price list" sheet
#   - D column held "=E / 1.1" helper formulas for rows 16-22 and 26-32; these
#     were cleared out (left blank), while the shared-formula anchor for each
#     block shifts down to the next populated row (23 / 33) automatically.
#   - E16 got its literal rate value tweaked (rounded) from 306.075 to 306.08.
#   - The sheet becomes the active tab, scrolled back to the top of its
#     frozen pane, with D32 selected.
# ---------------------------------------------------------------------------

$stdClearRows = 16, 17, 18, 19, 20, 21, 22, 26, 27, 28, 29, 30, 31, 32
foreach ($r in $stdClearRows) {
    $wsStd.Range("D$r").ClearContents()
}

$wsStd.Range("E16").Value = 306.08

# ---------------------------------------------------------------------------
# "IAM UAW East" sheet
#   - D column held literal rate values for rows 16-22 and 26-32; these were
#     cleared out, and the corresponding E column (which held "=D*1.1")
#     formulas) now stores the literal, rounded result instead. The shared
#     formula anchor shifts down to the next populated row (23 / 33).
#   - The sheet is no longer the active tab; selection parked on E9, scrolled
#     back to the top of its frozen pane.
# ---------------------------------------------------------------------------

$uawClearRows = 16, 17, 18, 19, 20, 21, 22, 26, 27, 28, 30, 31, 32
foreach ($r in $uawClearRows) {
    $wsUaw.Range("D$r").ClearContents()
}

$uawLiteralE = @{
    16 = 218.63
    17 = 56.57
    18 = 64.5
    19 = 73.89
    20 = 80.41
    21 = 88.92
    22 = 102.28
    26 = 63.46
    27 = 76.09
    28 = 86.11
    29 = 0
    30 = 109.19
    31 = 134.76
    32 = 341
}
foreach ($r in $uawLiteralE.Keys) {
    $wsUaw.Range("E$r").Value = $uawLiteralE[$r]
}

# ---------------------------------------------------------------------------
# View / selection state: the standard price list becomes the active sheet
# (previously it was "IAM UAW East"); both sheets' panes scroll back to the
# top and get a fresh selected cell. Selecting a range activates its sheet,
# so select on "IAM UAW East" first and finish on the standard price list so
# that one is left as the active tab (matching workbookView's activeTab).
# ---------------------------------------------------------------------------

[void]$wsUaw.Range("E9").Select()
[void]$wsStd.Range("D32").Select()
